$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Sep 26 21:13:25 EDT 2023"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Sep 26 21:13:36 EDT 2023"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Tue Sep 26 21:13:46 EDT 2023"
$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Tue Sep 26 21:14:03 EDT 2023"
$ws.Range("A6").Value = "Fail"
$ws.Range("B6").Value = "Tue Sep 26 21:14:13 EDT 2023"
$ws.Range("A7").Value = "Pass"
$ws.Range("B7").Value = "Tue Sep 26 21:14:23 EDT 2023"
$ws.Range("A8").Value = "Pass"
$ws.Range("B8").Value = "Tue Sep 26 21:14:34 EDT 2023"
$ws.Range("A9").Value = "Pass"
$ws.Range("B9").Value = "Tue Sep 26 21:14:45 EDT 2023"
$ws.Range("A10").Value = "Pass"
$ws.Range("B10").Value = "Tue Sep 26 21:14:56 EDT 2023"
$ws.Range("A11").Value = "Fail"
$ws.Range("B11").Value = "Tue Sep 26 21:15:06 EDT 2023"
$ws.Range("A12").Value = "Fail"
$ws.Range("B12").Value = "Tue Sep 26 21:15:17 EDT 2023"
$ws.Range("A13").Value = "Pass"
$ws.Range("B13").Value = "Tue Sep 26 21:15:27 EDT 2023"
$ws.Range("A14").Value = "Pass"
$ws.Range("B14").Value = "Tue Sep 26 21:15:38 EDT 2023"
$ws.Range("A15").Value = "Fail"
$ws.Range("B15").Value = "Tue Sep 26 21:15:49 EDT 2023"
$ws.Range("A16").Value = "Fail"
$ws.Range("B16").Value = "Tue Sep 26 21:15:59 EDT 2023"
$ws.Range("A17").Value = "Fail"
$ws.Range("B17").Value = "Tue Sep 26 21:16:08 EDT 2023"
$ws.Range("A18").Value = "Fail"
$ws.Range("B18").Value = "Tue Sep 26 21:16:18 EDT 2023"
$ws.Range("A19").Value = "Pass"
$ws.Range("B19").Value = "Tue Sep 26 21:16:29 EDT 2023"
$ws.Range("A20").Value = "Pass"
$ws.Range("B20").Value = "Tue Sep 26 21:16:39 EDT 2023"
$ws.Range("A21").Value = "Fail"
$ws.Range("B21").Value = "Tue Sep 26 21:16:50 EDT 2023"
$ws.Range("A22").Value = "Fail"
$ws.Range("B22").Value = "Tue Sep 26 21:17:00 EDT 2023"
$ws.Range("A23").Value = "Fail"
$ws.Range("B23").Value = "Tue Sep 26 21:17:10 EDT 2023"
$ws.Range("A24").Value = "Fail"
$ws.Range("B24").Value = "Tue Sep 26 21:17:20 EDT 2023"
$ws.Range("A25").Value = "Pass"
$ws.Range("B25").Value = "Tue Sep 26 21:17:29 EDT 2023"
$ws.Range("A26").Value = "Fail"
$ws.Range("B26").Value = "Tue Sep 26 21:17:40 EDT 2023"
$ws.Range("A27").Value = "Fail"
$ws.Range("B27").Value = "Tue Sep 26 21:17:50 EDT 2023"
$ws.Range("A28").Value = "Fail"
$ws.Range("B28").Value = "Tue Sep 26 21:18:00 EDT 2023"
$ws.Range("A29").Value = "Fail"
$ws.Range("B29").Value = "Tue Sep 26 21:18:10 EDT 2023"
